$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the header row and first data row. Values are assigned in this
# particular order so the workbook's shared-string table is built up in
# the same sequence as the source data (category, sector, description,
# emissions, Fliegen, Name) rather than left-to-right cell order.
$ws.Range("B1").Value = "category"
$ws.Range("A1").Value = "sector"
$ws.Range("D1").Value = "description"
$ws.Range("E1").Value = "emissions [kgCO2/Pers/a]"
$ws.Range("A2").Value = "Fliegen"
$ws.Range("C1").Value = "Name"

# Widen the description column and select the "Name" header cell, matching
# the authoring session's final view state.
$ws.Columns.Item(4).ColumnWidth = 42
$ws.Range("C1").Select() | Out-Null
